# Updates cryptos list figures (price + 1h volume change) to the latest
# scraped values, as produced by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a "Price" (column D) cell while keeping it
# stored as text, even when the new value happens to look like a plain
# number (Excel would otherwise silently convert it to a numeric cell).
function Set-TextValue($cellRange, [string]$text) {
    if ($text -match '^[0-9]+(\.[0-9]+)?$') {
        $cellRange.NumberFormat = "@"
    }
    $cellRange.Value = $text
}

# Rows that only change the Volume(1h) column (E), price (D) stays the same.
$volumeOnly = @{
    4  = "  +0.23%  "
    7  = "  +3.92%  "
    8  = "  +0.08%  "
    9  = "  +3.79%  "
    11 = "  +3.81%  "
    13 = "  +2.17%  "
    19 = "  +10.63%  "
    20 = "  +3.84%  "
    22 = "  +2.15%  "
    25 = "  +7.03%  "
    26 = "  +0.30%  "
    29 = "  +0.74%  "
    30 = "  -0.16%  "
    40 = "  +0.95%  "
    43 = "  +2.84%  "
    45 = "  +0.07%  "
    48 = "  +1.27%  "
}

foreach ($row in $volumeOnly.Keys) {
    $ws.Range("E$row").Value = $volumeOnly[$row]
}

# Rows that change both Price (D) and Volume(1h) (E).
$priceAndVolume = @(
    @{ Row = 2;  D = "51.927.47"; E = "  +3.48%  " }
    @{ Row = 3;  D = "2.782.02";  E = "  +3.58%  " }
    @{ Row = 5;  D = "343.38";    E = "  +4.48%  " }
    @{ Row = 6;  D = "115.55";    E = "  +1.38%  " }
    @{ Row = 10; D = "42.49";     E = "  +5.66%  " }
    @{ Row = 12; D = "20.06";     E = "  -0.34%  " }
    @{ Row = 14; D = "7.64";      E = "  +0.40%  " }
    @{ Row = 15; D = "3.215.73";  E = "  +4.30%  " }
    @{ Row = 16; D = "2.752.48";  E = "  +3.67%  " }
    @{ Row = 17; D = "0.883";     E = "  +1.25%  " }
    @{ Row = 18; D = "51.817.81"; E = "  +3.73%  " }
    @{ Row = 21; D = "13.29";     E = "  -3.00%  " }
    @{ Row = 23; D = "270.05";    E = "  -1.77%  " }
    @{ Row = 24; D = "70.05";     E = "  +0.42%  " }
    @{ Row = 27; D = "0.999";     E = "  -0.05%  " }
    @{ Row = 28; D = "10.23";     E = "  -0.48%  " }
    @{ Row = 31; D = "34.57";     E = "  -2.01%  " }
    @{ Row = 32; D = "50.14";     E = "  +0.90%  " }
    @{ Row = 33; D = "5.72";      E = "  +2.95%  " }
    @{ Row = 34; D = "0.0821";    E = "  -0.38%  " }
    @{ Row = 38; D = "2.10";      E = "  +1.18%  " }
    @{ Row = 39; D = "4.95";      E = "  -0.87%  " }
    @{ Row = 41; D = "2.66";      E = "  +23.87%  " }
    @{ Row = 42; D = "23.56";     E = "  -0.67%  " }
    @{ Row = 44; D = "126.48";    E = "  -1.29%  " }
    @{ Row = 46; D = "3.33";      E = "  -1.32%  " }
    @{ Row = 47; D = "2.065.35";  E = "  -0.55%  " }
    @{ Row = 49; D = "5.55";      E = "  +3.63%  " }
    @{ Row = 50; D = "0.901";     E = "  +12.93%  " }
    @{ Row = 51; D = "8.87";      E = "  -1.77%  " }
)

foreach ($item in $priceAndVolume) {
    Set-TextValue $ws.Range("D" + $item.Row) $item.D
    $ws.Range("E" + $item.Row).Value = $item.E
}

# Rows 35-37 got re-ranked: VeChain moved up to rank 35, FirstDigitalUSD
# shifted to rank 36, Celestia shifted to rank 37 (each row keeps its
# rank-index in column A, but coin/link/price/volume are replaced).
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D35") "0.0404"
$ws.Range("E35").Value = "  +14.65%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D36") "1.00"
$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D37") "19.05"
$ws.Range("E37").Value = "  -1.51%  "
